$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Widen column F to hold the explanation text
$ws.Columns.Item(6).ColumnWidth = 27.3

# --- Column A (labels) first, top to bottom ---
$ws.Range("A23").Value = "RF Hyperparameters"
$ws.Range("A23").Font.Bold = $true

$ws.Range("A24").Value = "n"
$ws.Range("A25").Value = "n_estimators"
$ws.Range("A26").Value = "max_depth"
$ws.Range("A27").Value = "min_samples_split"
$ws.Range("A28").Value = "min_samples_leaf"
$ws.Range("A29").Value = "max_features"

# --- Row 29 values (text, right aligned) before moving to the next label ---
$ws.Range("B29").Value = "log2"
$ws.Range("C29").Value = "sqrt"
$ws.Range("D29").Value = "sqrt"
$ws.Range("E29").Value = "log2"
$ws.Range("B29:E29").HorizontalAlignment = -4152

$ws.Range("A30").Value = "n_components"

# --- Columns B:E (numeric values) ---
$ws.Range("B24").Value = 9
$ws.Range("C24").Value = 12
$ws.Range("D24").Value = 15
$ws.Range("E24").Value = 18

$ws.Range("B25").Value = 500
$ws.Range("C25").Value = 500
$ws.Range("D25").Value = 500
$ws.Range("E25").Value = 500

$ws.Range("B26").Value = 20
$ws.Range("C26").Value = 25
$ws.Range("D26").Value = 29
$ws.Range("E26").Value = 17

$ws.Range("B27").Value = 4
$ws.Range("C27").Value = 2
$ws.Range("D27").Value = 2
$ws.Range("E27").Value = 14

$ws.Range("B28").Value = 1
$ws.Range("C28").Value = 1
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = 17

$ws.Range("B30").Value = 9
$ws.Range("C30").Value = 11
$ws.Range("D30").Value = 4
$ws.Range("E30").Value = 15

# --- Column F (explanations), top to bottom ---
$ws.Range("F24").Value = "Explanation"
$ws.Range("F25").Value = "number of trees per forest"
$ws.Range("F26").Value = "maximum depth of each tree"
$ws.Range("F27").Value = "minimum number of samples for each split"
$ws.Range("F28").Value = "minimum number of samples for each leaf node"
$ws.Range("F29").Value = "maximum number of features in each tree"
$ws.Range("F30").Value = "number of PCA components"

# Update selection to match the authored state
$ws.Range("A23").Select()
